# Insert a new weekly price record row at row 104 ("Fruta / hortaliza, semanal").
# This pushes the existing rows 104-197 down to 105-198 and fills the new
# row with the latest observation for "Arveja Verde".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(104).Insert()

$ws.Range("A104").Value = 9
$ws.Range("B104").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C104").Value = "Metropolitana"
$ws.Range("D104").Value = 45240
$ws.Range("E104").Value = 13
$ws.Range("F104").Value = 100112022
$ws.Range("G104").Value = "Arveja Verde"
$ws.Range("H104").Value = "Sin especificar"
$ws.Range("I104").Value = "Primera"
$ws.Range("J104").Value = 52
$ws.Range("K104").Value = 16000
$ws.Range("L104").Value = 17000
$ws.Range("M104").Value = 16500
$ws.Range("N104").Value = "`$/saco 25 kilos"
$ws.Range("O104").Value = "Región del Maule"
$ws.Range("P104").Value = 660
$ws.Range("Q104").Value = 25
$ws.Range("R104").Value = "Hortaliza"
